$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "M3"
$ws.Range("H2").Value = "M1"
$ws.Range("I2").Value = "A1"
$ws.Range("P2").Value = "A1"
$ws.Range("Q2").Value = "A1"
$ws.Range("R2").Value = "M1"
$ws.Range("S2").Value = "M1"
$ws.Range("T2").Value = "DO"
$ws.Range("U2").Value = "M1"
$ws.Range("V2").Value = "M3"
$ws.Range("W2").Value = "DO"
$ws.Range("Y2").Value = "PH"
$ws.Range("Z2").Value = "PH"
$ws.Range("B3").Value = "A1"
$ws.Range("C3").Value = "A2"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = "A2"
$ws.Range("F3").Value = "DO"
$ws.Range("G3").Value = "M2"
$ws.Range("I3").Value = "DO"
$ws.Range("J3").Value = "M2"
$ws.Range("K3").Value = "A1"
$ws.Range("M3").Value = "M2"
$ws.Range("N3").Value = "M2"
$ws.Range("O3").Value = "A2"
$ws.Range("R3").Value = "A2"
$ws.Range("W3").Value = "M2"
$ws.Range("X3").Value = "A2"
$ws.Range("Y3").Value = "PH"
$ws.Range("Z3").Value = "PH"
$ws.Range("AA3").Value = "DO"
$ws.Range("AB3").Value = "M2"
$ws.Range("C4").Value = "M1"
$ws.Range("D4").Value = "M3"
$ws.Range("E4").Value = "M1"
$ws.Range("F4").Value = "A1"
$ws.Range("G4").Value = "A1"
$ws.Range("H4").Value = "DO"
$ws.Range("I4").Value = "M1"
$ws.Range("J4").Value = "A1"
$ws.Range("K4").Value = "DO"
$ws.Range("L4").Value = "M3"
$ws.Range("M4").Value = "A1"
$ws.Range("O4").Value = "M1"
$ws.Range("P4").Value = "DO"
$ws.Range("Q4").Value = "M1"
$ws.Range("R4").Value = "A1"
$ws.Range("S4").Value = "A1"
$ws.Range("T4").Value = "M3"
$ws.Range("Y4").Value = "PH"
$ws.Range("Z4").Value = "PH"
$ws.Range("AB4").Value = "A1"
$ws.Range("AC4").Value = "M1"
$ws.Range("E5").Value = "M1"
$ws.Range("F5").Value = "A2"
$ws.Range("G5").Value = "M1"
$ws.Range("I5").Value = "DO"
$ws.Range("J5").Value = "M2"
$ws.Range("K5").Value = "A2"
$ws.Range("L5").Value = "A2"
$ws.Range("M5").Value = "M2"
$ws.Range("N5").Value = "M1"
$ws.Range("O5").Value = "M1"
$ws.Range("P5").Value = "DO"
$ws.Range("Q5").Value = "M2"
$ws.Range("S5").Value = "M2"
$ws.Range("U5").Value = "M1"
$ws.Range("W5").Value = "A2"
$ws.Range("X5").Value = "DO"
$ws.Range("Y5").Value = "PH"
$ws.Range("Z5").Value = "PH"
$ws.Range("AB5").Value = "M2"
$ws.Range("AC5").Value = "A2"
$ws.Range("C6").Value = "M2"
$ws.Range("E6").Value = "A2"
$ws.Range("F6").Value = "M1"
$ws.Range("H6").Value = "A1"
$ws.Range("I6").Value = "M2"
$ws.Range("K6").Value = "A2"
$ws.Range("L6").Value = "A2"
$ws.Range("M6").Value = "M1"
$ws.Range("N6").Value = "M1"
$ws.Range("O6").Value = "DO"
$ws.Range("P6").Value = "M1"
$ws.Range("Q6").Value = "M1"
$ws.Range("R6").Value = "A1"
$ws.Range("S6").Value = "A1"
$ws.Range("T6").Value = "M1"
$ws.Range("U6").Value = "M3"
$ws.Range("V6").Value = "DO"
$ws.Range("W6").Value = "A2"
$ws.Range("Y6").Value = "PH"
$ws.Range("Z6").Value = "PH"
$ws.Range("AA6").Value = "DO"
$ws.Range("AB6").Value = "M2"
$ws.Range("D7").Value = "M1"
$ws.Range("E7").Value = "DO"
$ws.Range("F7").Value = "A1"
$ws.Range("H7").Value = "M3"
$ws.Range("K7").Value = "M3"
$ws.Range("L7").Value = "DO"
$ws.Range("O7").Value = "A1"
$ws.Range("R7").Value = "M3"
$ws.Range("S7").Value = "DO"
$ws.Range("T7").Value = "A1"
$ws.Range("X7").Value = "M3"
$ws.Range("Y7").Value = "PH"
$ws.Range("Z7").Value = "PH"
$ws.Range("AC7").Value = "DO"
$ws.Range("C8").Value = "M2"
$ws.Range("D8").Value = "A2"
$ws.Range("E8").Value = "M1"
$ws.Range("F8").Value = "A1"
$ws.Range("H8").Value = "DO"
$ws.Range("I8").Value = "A2"
$ws.Range("J8").Value = "A1"
$ws.Range("K8").Value = "M1"
$ws.Range("L8").Value = "A2"
$ws.Range("O8").Value = "DO"
$ws.Range("P8").Value = "M2"
$ws.Range("R8").Value = "A2"
$ws.Range("S8").Value = "M2"
$ws.Range("T8").Value = "A1"
$ws.Range("V8").Value = "DO"
$ws.Range("Y8").Value = "PH"
$ws.Range("Z8").Value = "PH"
$ws.Range("AA8").Value = "DO"
$ws.Range("AB8").Value = "M2"
$ws.Range("AC8").Value = "A2"
$ws.Range("B9").Value = "DO"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = "A2"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "M1"
$ws.Range("H9").Value = "A1"
$ws.Range("I9").Value = "DO"
$ws.Range("J9").Value = "M2"
$ws.Range("K9").Value = "A2"
$ws.Range("L9").Value = "M2"
$ws.Range("M9").Value = "M1"
$ws.Range("N9").Value = "A2"
$ws.Range("O9").Value = "M1"
$ws.Range("P9").Value = "M2"
$ws.Range("R9").Value = "A2"
$ws.Range("T9").Value = "M1"
$ws.Range("U9").Value = "A2"
$ws.Range("W9").Value = "A2"
$ws.Range("X9").Value = "M2"
$ws.Range("Y9").Value = "PH"
$ws.Range("Z9").Value = "PH"
$ws.Range("AB9").Value = "M2"
$ws.Range("AC9").Value = "DO"
$ws.Range("B10").Value = "DO"
$ws.Range("C10").Value = "A2"
$ws.Range("D10").Value = "M1"
$ws.Range("F10").Value = "M2"
$ws.Range("G10").Value = "A1"
$ws.Range("H10").Value = "M2"
$ws.Range("I10").Value = "M1"
$ws.Range("J10").Value = "DO"
$ws.Range("K10").Value = "M2"
$ws.Range("L10").Value = "M2"
$ws.Range("M10").Value = "A2"
$ws.Range("N10").Value = "M1"
$ws.Range("O10").Value = "A2"
$ws.Range("Q10").Value = "A1"
$ws.Range("R10").Value = "M1"
$ws.Range("S10").Value = "A2"
$ws.Range("T10").Value = "DO"
$ws.Range("U10").Value = "M2"
$ws.Range("V10").Value = "A2"
$ws.Range("X10").Value = "A2"
$ws.Range("Y10").Value = "PH"
$ws.Range("Z10").Value = "PH"
